$wb = $excel.ActiveWorkbook

# "year_Vecteurs" sheet: the B column ("year") values for the second block of
# rows (rows 7-11) move from 2020 to 2050.
$wsYearVec = $wb.Worksheets.Item("year_Vecteurs")
$wsYearVec.Range("B7").Value = 2050
$wsYearVec.Range("B8").Value = 2050
$wsYearVec.Range("B9").Value = 2050
$wsYearVec.Range("B10").Value = 2050
$wsYearVec.Range("B11").Value = 2050

# Make "year_Vecteurs" the active sheet/tab (it was "retrofit_Transition"
# before) and leave the selection on G19, matching the saved view state.
$wsYearVec.Activate()
$wsYearVec.Range("G19").Select()
